# Commit Final v2 - Agregado el link del repositorio a la presentacion
#
# The last slide ("Videos de la Presentación y Demostración" / YouTube
# playlist link) is duplicated to the end of the deck, unchanged. The
# original slide is then repurposed to point at the GitHub repository
# that hosts this presentation.

$p = $ppt.ActivePresentation

# Last slide in the (7-slide) deck - the "Videos..." slide.
$lastIndex = $p.Slides.Count
$original = $p.Slides.Item($lastIndex)

# Duplicate it; the duplicate is appended right after it (becomes the
# new last slide) and keeps the old "Videos..." / YouTube content.
$duplicate = $original.Duplicate()

# Re-purpose the original slide: title becomes "Repositorio de Github"
# (as two runs, since "Github" is typed/flagged separately), and the
# body becomes the GitHub repository URL.
$titleRange = $original.Shapes.Item(1).TextFrame.TextRange
$titleRange.Text = "Repositorio de "
$titleRange.InsertAfter("Github") | Out-Null

$bodyRange = $original.Shapes.Item(2).TextFrame.TextRange
$bodyRange.Text = "https://github.com/ZiraelS/Presentacion_Pivot_Unpivot_Lookup"
